# Edit script for "Preparing for a staging release.docx"
#
# 1. Remove the _GoBack bookmark from its old location (end of the
#    paragraph ending in "...delay-signing).").
# 2. Append new explanatory text (several runs) to the end of the
#    paragraph "Check the CloudApiPublic and Sample-Live-Sync references ".
# 3. Re-add the _GoBack bookmark at its new location (end of the
#    paragraph ending in "...file, always answer \u201cNO\u201d.").

$d = $word.ActiveDocument

# --- Step 1: remove old _GoBack bookmark -----------------------------
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
    # bookmark not present - nothing to do
}

# --- Step 2: insert the new explanatory runs --------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "Check the CloudApiPublic and Sample-Live-Sync references ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)

    $rng.InsertAfter("(")
    $rng.Collapse(0)

    $rng.InsertAfter("CloudApiPublic")
    $rng.Collapse(0)

    $rng.InsertAfter(" should appear with no warning symbol for the ")
    $rng.Collapse(0)

    $rng.InsertAfter("BadgeCOMLib")
    $rng.Collapse(0)

    $rng.InsertAfter(" reference, but Sample-Live-Sync should actually have a warning symbol on the Cloud reference since its correct version requirement should be for the new version of ")
    $rng.Collapse(0)

    $rng.InsertAfter("CloudApiPublic")
    $rng.Collapse(0)

    $rng.InsertAfter(" which hasn" + [char]0x2019 + "t been built yet)")
    $rng.Collapse(0)
}

# --- Step 3: add _GoBack bookmark at the new location ------------------
# The target position sits exactly on a paragraph boundary (the very end
# of the paragraph's text, right before the paragraph mark). Adding a
# zero-width bookmark straight at such a boundary position is handled
# unreliably by some engines, so nudge around it: temporarily insert a
# marker character after the target position (making the position an
# "interior" one), add the bookmark there, then remove the marker again.
# Because bookmarks are zero-width or collapse, they keep pointing at the
# same character offset once the marker is deleted.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "file, always answer " + [char]0x201C + "NO" + [char]0x201D + ".",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $rng2.Collapse(0)
    $bmPos = $rng2.Start

    $rng2.InsertAfter("X")

    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $markerRange = $d.Range($bmPos, $bmPos + 1)
    $markerRange.Delete()
}
